$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new labels and formulas in column H summarizing part counts
$ws.Range("H33").Value = "Unique parts (loaded) count:"
$ws.Range("H39").Value = "Number of parts:"
$ws.Range("H36").Value = "Number of loaded SMD parts:"

$ws.Range("H34").Formula = "=COUNTBLANK(I2:I30)"
$ws.Range("H37").Formula = "=SUMIF(I2:I30, """", A2:A30)"
$ws.Range("H40").Formula = "=SUM(A2:A30)"

# Update the view so the newly added column is visible
$ws.Range("H37").Select()
